$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.520808
$ws.Range("H2").Value = 7.562424
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7132103333333334
$ws.Range("N2").Value = 2.139631
$ws.Range("O2").Value = 0.3847801132923274
$ws.Range("P2").Value = 0.3847801132923274
$ws.Range("Q2").Value = 1.797866313949334
$ws.Range("R2").Value = 16.180796825544
$ws.Range("S2").Value = 0.3847801132923274
$ws.Range("T2").Value = 0.3847801132923274

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.520808
$ws.Range("H3").Value = 7.562424
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6105696666666667
$ws.Range("N3").Value = 1.831709
$ws.Range("O3").Value = 0.3294050219587282
$ws.Range("P3").Value = 0.3294050219587283
$ws.Range("Q3").Value = 1.539128900290667
$ws.Range("R3").Value = 13.852160102616
$ws.Range("S3").Value = 0.3294050219587282
$ws.Range("T3").Value = 0.3294050219587283

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.520808
$ws.Range("H4").Value = 7.562424
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.529773
$ws.Range("N4").Value = 1.589319
$ws.Range("O4").Value = 0.2858148647489443
$ws.Range("P4").Value = 0.2858148647489444
$ws.Range("Q4").Value = 1.335456016584
$ws.Range("R4").Value = 12.019104149256
$ws.Range("S4").Value = 0.2858148647489443
$ws.Range("T4").Value = 0.2858148647489444
